$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 12.92400781198218
$ws.Range("C2").Value = 9.209067715088828
$ws.Range("D2").Value = 5.972409321078159
$ws.Range("E2").Value = 11.1291814565961
$ws.Range("G2").Value = 3.614323427774079
$ws.Range("I2").Value = 18.8760028232632
$ws.Range("M2").Value = 14.95109763546927
$ws.Range("O2").Value = 20.62448068331636
$ws.Range("B3").Value = 12.28872745133855
$ws.Range("C3").Value = 8.642616897221853
$ws.Range("D3").Value = 5.851603696232399
$ws.Range("E3").Value = 11.04978727977199
$ws.Range("G3").Value = 3.616856179701941
$ws.Range("I3").Value = 18.99947721072358
$ws.Range("M3").Value = 14.64557723663231
$ws.Range("O3").Value = 20.66638700086532
$ws.Range("B4").Value = 11.88213560801137
$ws.Range("C4").Value = 8.273770014285866
$ws.Range("D4").Value = 5.777992244448566
$ws.Range("E4").Value = 11.00556831335975
$ws.Range("G4").Value = 3.618492507154332
$ws.Range("I4").Value = 19.08114008881916
$ws.Range("M4").Value = 14.45799390385102
$ws.Range("O4").Value = 20.6995217855982
$ws.Range("B5").Value = 11.71248220766868
$ws.Range("C5").Value = 8.118186301448187
$ws.Range("D5").Value = 5.748183176981053
$ws.Range("E5").Value = 10.98870241118656
$ws.Range("G5").Value = 3.619179814108593
$ws.Range("I5").Value = 19.1158828446871
$ws.Range("M5").Value = 14.38166228376551
$ws.Range("O5").Value = 20.71487583321625
$ws.Range("B6").Value = 11.68407821926929
$ws.Range("C6").Value = 8.092033719082727
$ws.Range("D6").Value = 5.743246136516417
$ws.Range("E6").Value = 10.98597189846408
$ws.Range("G6").Value = 3.619295180519137
$ws.Range("I6").Value = 19.12174012018251
$ws.Range("M6").Value = 14.36899733122284
$ws.Range("O6").Value = 20.71753685256045
$ws.Range("B7").Value = 11.87986337562923
$ws.Range("C7").Value = 8.271693092434564
$ws.Range("D7").Value = 5.777589405220151
$ws.Range("E7").Value = 11.00533616522208
$ws.Range("G7").Value = 3.618501693359864
$ws.Range("I7").Value = 19.08160272054494
$ws.Range("M7").Value = 14.45696387852172
$ws.Range("O7").Value = 20.69972137425093
$ws.Range("B8").Value = 12.70849896101933
$ws.Range("C8").Value = 9.01813061361263
$ws.Range("D8").Value = 5.930666247564038
$ws.Range("E8").Value = 11.10087593151599
$ws.Range("G8").Value = 3.615179908197957
$ws.Range("I8").Value = 18.91735897274603
$ws.Range("M8").Value = 14.84581687633085
$ws.Range("O8").Value = 20.63738725629981
$ws.Range("B9").Value = 14.19537575866203
$ws.Range("C9").Value = 10.31454337039585
$ws.Range("D9").Value = 6.233316901555101
$ws.Range("E9").Value = 11.32339132069279
$ws.Range("G9").Value = 3.609307048853848
$ws.Range("I9").Value = 18.64197477409279
$ws.Range("M9").Value = 15.60379449885628
$ws.Range("O9").Value = 20.574319628567
$ws.Range("B10").Value = 15.19610148041843
$ws.Range("C10").Value = 11.16485671810286
$ws.Range("D10").Value = 6.4545756830684
$ws.Range("E10").Value = 11.5070981103753
$ws.Range("G10").Value = 3.605378647054978
$ws.Range("I10").Value = 18.46851750094643
$ws.Range("M10").Value = 16.15196158740562
$ws.Range("O10").Value = 20.56455825917545
$ws.Range("B11").Value = 15.63029075770203
$ws.Range("C11").Value = 11.5295481698939
$ws.Range("D11").Value = 6.554474854435558
$ws.Range("E11").Value = 11.59475967335311
$ws.Range("G11").Value = 3.603674461764753
$ws.Range("I11").Value = 18.3959641025953
$ws.Range("M11").Value = 16.39827419626991
$ws.Range("O11").Value = 20.56814224800646
$ws.Range("B12").Value = 15.79160637516717
$ws.Range("C12").Value = 11.66447262277477
$ws.Range("D12").Value = 6.592155417966795
$ws.Range("E12").Value = 11.62851524534095
$ws.Range("G12").Value = 3.603040974251074
$ws.Range("I12").Value = 18.3694112240755
$ws.Range("M12").Value = 16.49101701964127
$ws.Range("O12").Value = 20.5706580653702
$ws.Range("B13").Value = 15.75700321416837
$ws.Range("C13").Value = 11.63555539980852
$ws.Range("D13").Value = 6.584047491031018
$ws.Range("E13").Value = 11.6212209392178
$ws.Range("G13").Value = 3.603176881104533
$ws.Range("I13").Value = 18.37508874890632
$ws.Range("M13").Value = 16.47106818323029
$ws.Range("O13").Value = 20.57006464724983
$ws.Range("B14").Value = 15.64362485361723
$ws.Range("C14").Value = 11.54071210998454
$ws.Range("D14").Value = 6.55757807080027
$ws.Range("E14").Value = 11.59752570860727
$ws.Range("G14").Value = 3.603622107237737
$ws.Range("I14").Value = 18.39376106674275
$ws.Range("M14").Value = 16.40591533951033
$ws.Range("O14").Value = 20.56832598264228
$ws.Range("B15").Value = 15.57377128998064
$ws.Range("C15").Value = 11.48220438069958
$ws.Range("D15").Value = 6.541344197963964
$ws.Range("E15").Value = 11.5830837603456
$ws.Range("G15").Value = 3.603896362448393
$ws.Range("I15").Value = 18.40531864741892
$ws.Range("M15").Value = 16.36593564126185
$ws.Range("O15").Value = 20.5674120068473
$ws.Range("B16").Value = 15.16729585993089
$ws.Range("C16").Value = 11.14057850205355
$ws.Range("D16").Value = 6.448028348087803
$ws.Range("E16").Value = 11.50144917379324
$ws.Range("G16").Value = 3.605491681421161
$ws.Range("I16").Value = 18.47338755400342
$ws.Range("M16").Value = 16.13579511826643
$ws.Range("O16").Value = 20.56448594468502
$ws.Range("B17").Value = 14.91248719798207
$ws.Range("C17").Value = 10.92533914257917
$ws.Range("D17").Value = 6.390559549138374
$ws.Range("E17").Value = 11.45239741466743
$ws.Range("G17").Value = 3.606491534888371
$ws.Range("I17").Value = 18.51677843200608
$ws.Range("M17").Value = 15.99376288658857
$ws.Range("O17").Value = 20.56474994778901
$ws.Range("B18").Value = 14.76395059956845
$ws.Range("C18").Value = 10.799456661368
$ws.Range("D18").Value = 6.35743648598347
$ws.Range("E18").Value = 11.42457114595905
$ws.Range("G18").Value = 3.607074427624294
$ws.Range("I18").Value = 18.5423332256919
$ws.Range("M18").Value = 15.91178763269855
$ws.Range("O18").Value = 20.56565689816472
$ws.Range("B19").Value = 14.71332159700275
$ws.Range("C19").Value = 10.75647723963419
$ws.Range("D19").Value = 6.346211098464999
$ws.Range("E19").Value = 11.41521695995491
$ws.Range("G19").Value = 3.607273127339825
$ws.Range("I19").Value = 18.55108802277895
$ws.Range("M19").Value = 15.88398665713818
$ws.Range("O19").Value = 20.56609349990534
$ws.Range("B20").Value = 14.93981726457998
$ws.Range("C20").Value = 10.94846724324787
$ws.Range("D20").Value = 6.39668459684184
$ws.Range("E20").Value = 11.45757920584186
$ws.Range("G20").Value = 3.606384291624591
$ws.Range("I20").Value = 18.51209750196388
$ws.Range("M20").Value = 16.00891233481918
$ws.Range("O20").Value = 20.56464365537295
$ws.Range("B21").Value = 15.67701156190043
$ws.Range("C21").Value = 11.56865604028367
$ws.Range("D21").Value = 6.565357141471956
$ws.Range("E21").Value = 11.60447060675373
$ws.Range("G21").Value = 3.603491012430038
$ws.Range("I21").Value = 18.38825148067732
$ws.Range("M21").Value = 16.4250674069672
$ws.Range("O21").Value = 20.56880519440427
$ws.Range("B22").Value = 16.14070151663196
$ws.Range("C22").Value = 11.95547489536272
$ws.Range("D22").Value = 6.67470812430246
$ws.Range("E22").Value = 11.70372281594868
$ws.Range("G22").Value = 3.601669130795175
$ws.Range("I22").Value = 18.31268604130762
$ws.Range("M22").Value = 16.69391860569376
$ws.Range("O22").Value = 20.57827967745851
$ws.Range("B23").Value = 15.89490002051421
$ws.Range("C23").Value = 11.75071408203348
$ws.Range("D23").Value = 6.616439438322121
$ws.Range("E23").Value = 11.65046241010038
$ws.Range("G23").Value = 3.602635207302772
$ws.Range("I23").Value = 18.35252225951848
$ws.Range("M23").Value = 16.5507428967456
$ws.Range("O23").Value = 20.57260369876723
$ws.Range("B24").Value = 14.92746770126034
$ws.Range("C24").Value = 10.93801769035694
$ws.Range("D24").Value = 6.39391571971069
$ws.Range("E24").Value = 11.45523535019036
$ws.Range("G24").Value = 3.606432751180092
$ws.Range("I24").Value = 18.51421185455625
$ws.Range("M24").Value = 16.00206425595769
$ws.Range("O24").Value = 20.56468935814307
$ws.Range("B25").Value = 13.8087812832163
$ws.Range("C25").Value = 9.981775115710084
$ws.Range("D25").Value = 6.151452347207816
$ws.Range("E25").Value = 11.25954676489991
$ws.Range("G25").Value = 3.610827637337112
$ws.Range("I25").Value = 18.7114334889119
$ws.Range("M25").Value = 15.39985918335429
$ws.Range("O25").Value = 20.58499004382831
